$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the last row of data (row 6) with new values
$ws.Range("E6").Value = 7
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 13

# Match the saved cursor/selection position
$ws.Range("E6").Select()
